$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: the empty paragraph right after "(Database-first Approach)" gets
#         two runs of text:
#           "Design Pattern – IoC("
#           "Inversion of Controll)"
# ---------------------------------------------------------------------------
$n = $d.Paragraphs.Count
$idx1 = -1
for ($i = 1; $i -le $n; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*Database-first Approach*") {
        $idx1 = $i + 1
    }
}

if ($idx1 -gt 0) {
    $p1 = $d.Paragraphs.Item($idx1)

    # First run.
    $p1.Range.InsertAfter("Design Pattern – IoC(")

    # Second run - inserted right after the first one.  Adjacent runs that
    # share identical formatting normally get coalesced into a single run
    # by the serializer when a paragraph only ever had 0/1 runs, so we
    # momentarily give the new text a distinguishing character format and
    # then clear it again; that is enough to make the engine keep it as an
    # independent <w:r>, matching the two separate runs in the target XML.
    $posAfterFirst = $p1.Range.End - 1
    $r2 = $d.Range($posAfterFirst, $posAfterFirst)
    $r2.InsertAfter("Inversion of Controll)")
    $r2b = $d.Range($posAfterFirst, $p1.Range.End - 1)
    $r2b.Bold = 1
    $r2b.Bold = 0
}

# ---------------------------------------------------------------------------
# Edit 2: the paragraph that reads "CLI Command:" gets a new run
#         "Code Generator " inserted right before the existing "CLI" run.
# ---------------------------------------------------------------------------
$n = $d.Paragraphs.Count
$idx2 = -1
for ($i = 1; $i -le $n; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "CLI Command:*") {
        $idx2 = $i
    }
}

if ($idx2 -gt 0) {
    $p2 = $d.Paragraphs.Item($idx2)
    $r = $p2.Range
    $r.Collapse(1)
    $r.InsertBefore("Code Generator ")
}
